# Implement tri-training except for BERT, ElMo and CNN
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Resize the picture (Content Placeholder 4) ---
$pic = $s.Shapes.Item(1)
$pic.Height = 462.7573

# --- Resize the existing notes textbox (TextBox 6) and append a new bullet ---
$notes = $s.Shapes.Item(2)
$notes.Height = 421.67812

$notesRange = $notes.TextFrame.TextRange
$notesRange.InsertAfter("`rCompare the differences in performances with (possibly) self-training or co-training models") | Out-Null

# --- Add a new caption textbox below the picture, cloned from the notes box
#     so it inherits identical run/paragraph formatting ---
$caption = $notes.Duplicate().Item(1)
$caption.Name = "TextBox 1"
$caption.TextFrame.TextRange.Text = "Results in the red box have been reproduced"
$caption.Left = 58.7185
$caption.Top = 496.31072
$caption.Width = 446.67962
$caption.Height = 29.35922
